# Updated cryptos list on Sun Apr 28 08:06:56 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "601.04") are not silently coerced into floating point numbers.
    $c.NumberFormat = "@"
    $c.Value = $val
    # Remove the temporary formatting again so the cell keeps its original
    # (unstyled) look, matching the source workbook.
    $c.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "64.022.86"
Set-TextValue "E2" "  +1.51%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.319.26"
Set-TextValue "E3" "  +6.27%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "601.04"
Set-TextValue "E5" "  +1.40%  "

# Row 6 - Solana
Set-TextValue "D6" "143.66"
Set-TextValue "E6" "  +5.43%  "

# Row 7 - USDC
Set-TextValue "E7" "  -0.07%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.321.55"
Set-TextValue "E8" "  +6.58%  "

# Row 9 - XRP
Set-TextValue "E9" "  +1.56%  "

# Row 11 - Toncoin
Set-TextValue "D11" "5.55"
Set-TextValue "E11" "  +5.47%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.476"
Set-TextValue "E12" "  +4.46%  "

# Row 13 - ShibaInu
Set-TextValue "E13" "  +1.66%  "

# Row 14 - Avalanche
Set-TextValue "D14" "35.02"
Set-TextValue "E14" "  +2.74%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.864.34"
Set-TextValue "E15" "  +6.26%  "

# Row 16 - TRON
Set-TextValue "E16" "  +0.39%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "3.319.71"
Set-TextValue "E17" "  +6.27%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "64.096.97"
Set-TextValue "E18" "  +1.72%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.93"
Set-TextValue "E19" "  +3.97%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "483.30"
Set-TextValue "E20" "  +2.40%  "

# Row 21 - Chainlink
Set-TextValue "E21" "  +1.30%  "

# Row 22 - Polygon
Set-TextValue "E22" "  +6.33%  "

# Row 23 - Uniswap
Set-TextValue "D23" "8.02"
Set-TextValue "E23" "  +3.98%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "13.67"
Set-TextValue "E24" "  +5.72%  "

# Row 25 - Litecoin
Set-TextValue "D25" "84.92"
Set-TextValue "E25" "  -1.79%  "

# Row 26 - Dai
Set-TextValue "E26" "  +0.15%  "

# Row 27 - PancakeSwap
Set-TextValue "E27" "  +2.42%  "

# Row 28 - NEARProtocol
Set-TextValue "D28" "7.34"
Set-TextValue "E28" "  +4.99%  "

# Row 29 - FirstDigitalUSD
Set-TextValue "E29" "  -0.06%  "

# Row 30 - RenderToken
Set-TextValue "D30" "8.26"
Set-TextValue "E30" "  +4.22%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "29.75"
Set-TextValue "E31" "  +11.66%  "

# Row 32 - ImmutableX
Set-TextValue "E32" "  +5.60%  "

# Row 33 - Hedera
Set-TextValue "E33" "  -1.99%  "

# Row 34 - Stacks
Set-TextValue "E34" "  +2.58%  "

# Row 35 - Mantle
Set-TextValue "E35" "  +2.64%  "

# Row 36 - Filecoin
Set-TextValue "D36" "6.02"
Set-TextValue "E36" "  +3.85%  "

# Row 37 - PEPE
Set-TextValue "E37" "  +7.68%  "

# Row 38 - OKB
Set-TextValue "D38" "53.31"
Set-TextValue "E38" "  +2.32%  "

# Row 39 - VeChain
Set-TextValue "E39" "  +4.79%  "

# Row 40 - Bittensor
Set-TextValue "D40" "436.18"
Set-TextValue "E40" "  +2.87%  "

# Row 41 - Maker
Set-TextValue "D41" "3.069.63"
Set-TextValue "E41" "  +5.88%  "

# Rows 42 & 43 - ranking swapped: dogwifhat/Cosmos -> Cosmos/dogwifhat
Set-TextValue "B42" "Cosmos"
Set-TextValue "C42" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D42" "8.47"
Set-TextValue "E42" "  +3.17%  "

Set-TextValue "B43" "dogwifhat"
Set-TextValue "C43" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.78"
Set-TextValue "E43" "  +3.49%  "

# Row 44 - Kaspa
Set-TextValue "E44" "  +0.43%  "

# Row 45 - TheGraph
Set-TextValue "E45" "  +2.43%  "

# Row 46 - Fetch.AI
Set-TextValue "E46" "  +5.13%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "26.62"
Set-TextValue "E47" "  +4.39%  "

# Row 48 - Arweave
Set-TextValue "D48" "36.18"
Set-TextValue "E48" "  +15.53%  "

# Row 50 - Stellar
Set-TextValue "E50" "  +3.09%  "

# Row 51 - ThetaToken
Set-TextValue "E51" "  +2.05%  "
